$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.617.92'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '''1.920.78'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''245.78'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '''0.4851'
$ws.Range("E7").Value = '  +2.70%  '
$ws.Range("D8").Value = '''0.2910'
$ws.Range("E8").Value = '  +1.01%  '
$ws.Range("D9").Value = '''0.06721'
$ws.Range("E9").Value = '  -1.23%  '
$ws.Range("D10").Value = '''111.66'
$ws.Range("E10").Value = '  +6.20%  '
$ws.Range("D11").Value = '''19.33'
$ws.Range("E11").Value = '  +5.22%  '
$ws.Range("D12").Value = '''1.920.95'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").Value = '''0.07593'
$ws.Range("E13").Value = '  -1.27%  '
$ws.Range("D14").Value = '''5.340'
$ws.Range("E14").Value = '  +0.95%  '
$ws.Range("D15").Value = '''0.6731'
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").Value = '''294.48'
$ws.Range("E16").Value = '  +2.02%  '
$ws.Range("D17").Value = '''30.626.86'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").Value = '''13.07'
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '''0.000007562'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '''2.172.90'
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '''5.524'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = '''1.000'
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '''6.430'
$ws.Range("E24").Value = '  +1.96%  '
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("D26").Value = '''164.75'
$ws.Range("E26").Value = '  -2.47%  '
$ws.Range("D27").Value = '''20.30'
$ws.Range("E27").Value = '  -4.82%  '
$ws.Range("D28").Value = '''2.099'
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("E29").Value = '  +0.51%  '
$ws.Range("D30").Value = '''1.442'
$ws.Range("E30").Value = '  +3.32%  '
$ws.Range("D31").Value = '''4.132'
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("D32").Value = '''4.085'
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").Value = '''0.05029'
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("D34").Value = '''0.7407'
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("D35").Value = '''1.140'
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("D36").Value = '''0.9998'
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.02025'
$ws.Range("E37").Value = '  -2.12%  '
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").Value = '''2.700'
$ws.Range("E38").Value = '  -1.82%  '
$ws.Range("D39").Value = '''2.702'
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("D40").Value = '''110.00'
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("D41").Value = '''2.019'
$ws.Range("E41").Value = '  -1.98%  '
$ws.Range("D42").Value = '''0.4444'
$ws.Range("E42").Value = '  +0.86%  '
$ws.Range("D43").Value = '''0.8660'
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("D44").Value = '''5.883'
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").Value = '''70.27'
$ws.Range("E45").Value = '  +4.54%  '
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").Value = '''7.270'
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("D48").Value = '''48.52'
$ws.Range("E48").Value = '  -0.28%  '
$ws.Range("D49").Value = '''9.287'
$ws.Range("E49").Value = '  -0.55%  '
$ws.Range("D50").Value = '''0.1227'
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("D51").Value = '''0.2524'
$ws.Range("E51").Value = '  +2.74%  '
